# ---------------------------------------------------------------------------
# Replays, via PowerPoint COM-interop, the two changes recorded in the
# commit's OOXML diff:
#
#   1. The table on slide 5 gets a new table style applied
#      (tableStyleId {00EFBCB0-FCEC-4BED-84DC-87DE546FE399}
#                 -> {16C273CF-9C07-4F40-BE37-45F87A16DAE6}).
#
#   2. The deck's theme (ppt/theme/theme1.xml, the theme attached to the
#      slide master / used by every slide) is switched from the custom
#      "Integral" / "Red Violet" colour scheme to the stock Office theme
#      colour scheme. Concretely, every one of the twelve theme colour
#      slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) is reset to
#      the standard Office values; font scheme / format scheme were
#      already identical between the two themes, so only the colours
#      actually move.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 ---------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{16C273CF-9C07-4F40-BE37-45F87A16DAE6}")

# --- 2. Apply the stock "Office" colour scheme to the deck theme ------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Index : slot      : target RGB (hex)   : COM RGB (0xBBGGRR)
# 1     : dk1        000000               0
# 2     : lt1         FFFFFF              16777215
# 3     : dk2         44546A              6968388
# 4     : lt2         E7E6E6              15132391
# 5     : accent1     5B9BD5              13998939
# 6     : accent2     ED7D31              3243501
# 7     : accent3     A5A5A5              10855845
# 8     : accent4     FFC000              49407
# 9     : accent5     4472C4              12874308
# 10    : accent6     70AD47              4697456
# 11    : hlink       0563C1              12673797
# 12    : folHlink    954F72              7491477

$colorScheme.Item(1).RGB  = 0
$colorScheme.Item(2).RGB  = 16777215
$colorScheme.Item(3).RGB  = 6968388
$colorScheme.Item(4).RGB  = 15132391
$colorScheme.Item(5).RGB  = 13998939
$colorScheme.Item(6).RGB  = 3243501
$colorScheme.Item(7).RGB  = 10855845
$colorScheme.Item(8).RGB  = 49407
$colorScheme.Item(9).RGB  = 12874308
$colorScheme.Item(10).RGB = 4697456
$colorScheme.Item(11).RGB = 12673797
$colorScheme.Item(12).RGB = 7491477
